$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 4: release date slipped a day later ---
# -> recolor from orange to green (matches the already-green rows, e.g. A7:D7)
$ws.Range("A4").Value = 45168
$ws.Range("A4:D4").Interior.Color = 5296274   # RGB(146,208,80) FF92D050 green

# (B5:B16 are "days left" countdown formulas driven by TODAY() - their cached
# results move on their own every time the sheet recalculates, no edit needed)

# --- Insert a new row for the Gyarados ex Premium Collection release ---
$ws.Rows("17").Insert()

$ws.Range("A17").Value = 45233
$ws.Range("B17").Formula = "=IF(A17 = ""TBA"", ""TBA"", IF(A17 - TODAY() > 0, A17 - TODAY(), ""out now""))"
$ws.Range("C17").Value = "Gyarados ex Premium Collection"
$ws.Range("D17").Value = "check type of reverse Magikarp"

# --- Leave the selection where the editor was last working ---
$ws.Range("G7").Select()
